$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.90715
$ws.Range("H2").Value = 2.72145
$ws.Range("I2").Value = 0.01717809939998381
$ws.Range("J2").Value = 0.01717809939998381
$ws.Range("M2").Value = 2.685464
$ws.Range("N2").Value = 8.056392000000001
$ws.Range("O2").Value = 0.06676031826184478
$ws.Range("P2").Value = 0.06676031826184478
$ws.Range("Q2").Value = 2.4361186676
$ws.Range("R2").Value = 21.9250680084
$ws.Range("S2").Value = 0.001146815383076524
$ws.Range("T2").Value = 0.001146815383076524

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.90715
$ws.Range("H3").Value = 2.72145
$ws.Range("I3").Value = 0.01717809939998381
$ws.Range("J3").Value = 0.01717809939998381
$ws.Range("O3").Value = 0.02342101692711854
$ws.Range("P3").Value = 0.02342101692711854
$ws.Range("Q3").Value = 0.85464506515
$ws.Range("R3").Value = 7.69180558635
$ws.Range("S3").Value = 0.0004023285568227457
$ws.Range("T3").Value = 0.0004023285568227457

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.90715
$ws.Range("H4").Value = 2.72145
$ws.Range("I4").Value = 0.01717809939998381
$ws.Range("J4").Value = 0.01717809939998381
$ws.Range("M4").Value = 35.399925
$ws.Range("N4").Value = 106.199775
$ws.Range("O4").Value = 0.8800379597140142
$ws.Range("P4").Value = 0.8800379597140142
$ws.Range("Q4").Value = 32.11304196375001
$ws.Range("R4").Value = 289.01737767375
$ws.Range("S4").Value = 0.01511737954772629
$ws.Range("T4").Value = 0.01511737954772628

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.90715
$ws.Range("H5").Value = 2.72145
$ws.Range("I5").Value = 0.01717809939998381
$ws.Range("J5").Value = 0.01717809939998381
$ws.Range("M5").Value = 1.197942333333333
$ws.Range("N5").Value = 3.593827
$ws.Range("O5").Value = 0.02978070509702244
$ws.Range("P5").Value = 0.02978070509702244
$ws.Range("Q5").Value = 1.086713387683333
$ws.Range("R5").Value = 9.78042048915
$ws.Range("S5").Value = 0.0005115759123582561
$ws.Range("T5").Value = 0.000511575912358256

# Row 6
$ws.Range("I6").Value = 0.0237180037344858
$ws.Range("J6").Value = 0.0237180037344858
$ws.Range("M6").Value = 2.685464
$ws.Range("N6").Value = 8.056392000000001
$ws.Range("O6").Value = 0.06676031826184478
$ws.Range("P6").Value = 0.06676031826184478
$ws.Range("Q6").Value = 3.363577675877334
$ws.Range("R6").Value = 30.27219908289601
$ws.Range("S6").Value = 0.001583421477849895
$ws.Range("T6").Value = 0.001583421477849895

# Row 7
$ws.Range("I7").Value = 0.0237180037344858
$ws.Range("J7").Value = 0.0237180037344858
$ws.Range("O7").Value = 0.02342101692711854
$ws.Range("P7").Value = 0.02342101692711854
$ws.Range("S7").Value = 0.0005554997669428526
$ws.Range("T7").Value = 0.0005554997669428526

# Row 8
$ws.Range("I8").Value = 0.0237180037344858
$ws.Range("J8").Value = 0.0237180037344858
$ws.Range("M8").Value = 35.399925
$ws.Range("N8").Value = 106.199775
$ws.Range("O8").Value = 0.8800379597140142
$ws.Range("P8").Value = 0.8800379597140142
$ws.Range("Q8").Value = 44.33885446155001
$ws.Range("R8").Value = 399.0496901539501
$ws.Range("S8").Value = 0.02087274361498625
$ws.Range("T8").Value = 0.02087274361498625

# Row 9
$ws.Range("I9").Value = 0.0237180037344858
$ws.Range("J9").Value = 0.0237180037344858
$ws.Range("M9").Value = 1.197942333333333
$ws.Range("N9").Value = 3.593827
$ws.Range("O9").Value = 0.02978070509702244
$ws.Range("P9").Value = 0.02978070509702244
$ws.Range("Q9").Value = 1.500437946436223
$ws.Range("R9").Value = 13.503941517926
$ws.Range("S9").Value = 0.0007063388747067986
$ws.Range("T9").Value = 0.0007063388747067985

# Row 10
$ws.Range("G10").Value = 2.247832333333333
$ws.Range("H10").Value = 6.743497
$ws.Range("I10").Value = 0.04256571378106988
$ws.Range("J10").Value = 0.04256571378106987
$ws.Range("M10").Value = 2.685464
$ws.Range("N10").Value = 8.056392000000001
$ws.Range("O10").Value = 0.06676031826184478
$ws.Range("P10").Value = 0.06676031826184478
$ws.Range("Q10").Value = 6.036472809202666
$ws.Range("R10").Value = 54.328255282824
$ws.Range("S10").Value = 0.002841700599066818
$ws.Range("T10").Value = 0.002841700599066817

# Row 11
$ws.Range("G11").Value = 2.247832333333333
$ws.Range("H11").Value = 6.743497
$ws.Range("I11").Value = 0.04256571378106988
$ws.Range("J11").Value = 0.04256571378106987
$ws.Range("O11").Value = 0.02342101692711854
$ws.Range("P11").Value = 0.02342101692711854
$ws.Range("Q11").Value = 2.117730045712333
$ws.Range("R11").Value = 19.059570411411
$ws.Range("S11").Value = 0.0009969323029813207
$ws.Range("T11").Value = 0.0009969323029813204

# Row 12
$ws.Range("G12").Value = 2.247832333333333
$ws.Range("H12").Value = 6.743497
$ws.Range("I12").Value = 0.04256571378106988
$ws.Range("J12").Value = 0.04256571378106987
$ws.Range("M12").Value = 35.399925
$ws.Range("N12").Value = 106.199775
$ws.Range("O12").Value = 0.8800379597140142
$ws.Range("P12").Value = 0.8800379597140142
$ws.Range("Q12").Value = 79.57309601257501
$ws.Range("R12").Value = 716.1578641131749
$ws.Range("S12").Value = 0.03745944390966344
$ws.Range("T12").Value = 0.03745944390966343

# Row 13
$ws.Range("G13").Value = 2.247832333333333
$ws.Range("H13").Value = 6.743497
$ws.Range("I13").Value = 0.04256571378106988
$ws.Range("J13").Value = 0.04256571378106987
$ws.Range("M13").Value = 1.197942333333333
$ws.Range("N13").Value = 3.593827
$ws.Range("O13").Value = 0.02978070509702244
$ws.Range("P13").Value = 0.02978070509702244
$ws.Range("Q13").Value = 2.692773510335444
$ws.Range("R13").Value = 24.234961593019
$ws.Range("S13").Value = 0.001267636969358306
$ws.Range("T13").Value = 0.001267636969358306

# Row 14
$ws.Range("G14").Value = 48.40102466666667
$ws.Range("H14").Value = 145.203074
$ws.Range("I14").Value = 0.9165381830844606
$ws.Range("J14").Value = 0.9165381830844604
$ws.Range("M14").Value = 2.685464
$ws.Range("N14").Value = 8.056392000000001
$ws.Range("O14").Value = 0.06676031826184478
$ws.Range("P14").Value = 0.06676031826184478
$ws.Range("Q14").Value = 129.9792093054454
$ws.Range("R14").Value = 1169.812883749008
$ws.Range("S14").Value = 0.06118838080185154
$ws.Range("T14").Value = 0.06118838080185154

# Row 15
$ws.Range("G15").Value = 48.40102466666667
$ws.Range("H15").Value = 145.203074
$ws.Range("I15").Value = 0.9165381830844606
$ws.Range("J15").Value = 0.9165381830844604
$ws.Range("O15").Value = 0.02342101692711854
$ws.Range("P15").Value = 0.02342101692711854
$ws.Range("Q15").Value = 45.59962175998468
$ws.Range("R15").Value = 410.3965958398621
$ws.Range("S15").Value = 0.02146625630037162
$ws.Range("T15").Value = 0.02146625630037162

# Row 16
$ws.Range("G16").Value = 48.40102466666667
$ws.Range("H16").Value = 145.203074
$ws.Range("I16").Value = 0.9165381830844606
$ws.Range("J16").Value = 0.9165381830844604
$ws.Range("M16").Value = 35.399925
$ws.Range("N16").Value = 106.199775
$ws.Range("O16").Value = 0.8800379597140142
$ws.Range("P16").Value = 0.8800379597140142
$ws.Range("Q16").Value = 1713.39264312315
$ws.Range("R16").Value = 15420.53378810835
$ws.Range("S16").Value = 0.8065883926416383
$ws.Range("T16").Value = 0.8065883926416382

# Row 17
$ws.Range("G17").Value = 48.40102466666667
$ws.Range("H17").Value = 145.203074
$ws.Range("I17").Value = 0.9165381830844606
$ws.Range("J17").Value = 0.9165381830844604
$ws.Range("M17").Value = 1.197942333333333
$ws.Range("N17").Value = 3.593827
$ws.Range("O17").Value = 0.02978070509702244
$ws.Range("P17").Value = 0.02978070509702244
$ws.Range("Q17").Value = 57.9816364249109
$ws.Range("R17").Value = 521.8347278241981
$ws.Range("S17").Value = 0.02729515334059908
$ws.Range("T17").Value = 0.02729515334059908
